# aggiornamento fino a 20/09/2021
# Append 11 new daily rows (375-385) to the end of the data table,
# covering dates 2021-09-10 .. 2021-09-20 (serial 44449..44459).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(375, 44449, 1, 7, 69.51340615690168),
    @(376, 44450, 1, 4, 39.72194637537239),
    @(377, 44451, 2, 6, 59.5829195630586),
    @(378, 44452, 1, 5, 49.65243296921549),
    @(379, 44453, 0, 5, 49.65243296921549),
    @(380, 44454, 2, 7, 69.51340615690168),
    @(381, 44455, 0, 7, 69.51340615690168),
    @(382, 44456, 1, 7, 69.51340615690168),
    @(383, 44457, 0, 6, 59.5829195630586),
    @(384, 44458, 1, 5, 49.65243296921549),
    @(385, 44459, 0, 4, 39.72194637537239)
)

# Use the last existing data row (374, column A) as the format template so
# the newly appended date cells reuse the same style (bold, centered,
# bordered, custom date/time number format) instead of minting a new one.
$templateCell = $ws.Range("A374")
$templateCell.Copy() | Out-Null

foreach ($row in $data) {
    $r = $row[0]

    $destCell = $ws.Range("A$r")
    $destCell.PasteSpecial(-4122) | Out-Null

    $ws.Range("A$r").Value = $row[1]
    $ws.Range("B$r").Value = $row[2]
    $ws.Range("C$r").Value = $row[3]
    $ws.Range("D$r").Value = $row[4]
}

$excel.CutCopyMode = $false
